# Apply the "Clean Code Tip" update:
#   1. Refresh the cached text of every auto-updating
#      datetimeFigureOut date field (slide master + every slide layout)
#      from 5/19/2021 to 10/13/2021 - mirrors what PowerPoint does to
#      these fields whenever the deck is re-saved on a later date.
#   2. Update the cover slide's big title from
#      "Tip #1 / how to choose meaningful names"
#      to "Tip #6 / Avoid using too many arguments".

$p = $ppt.ActivePresentation

$ppPlaceholderDate = 16
$newDate = "10/13/2021"

# --- 1. Slide Master date placeholder -------------------------------
$master = $p.SlideMaster

function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $phType = $null
        try { $phType = $sh.PlaceholderFormat.Type } catch {}
        if ($phType -eq $ppPlaceholderDate -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Length -gt 0) {
                $tr.Characters(1, $tr.Length).Text = $newText
            }
        }
    }
}

Update-DatePlaceholder $master.Shapes $newDate

# --- 2. Every slide layout's date placeholder ------------------------
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes $newDate
}

# --- 3. Cover-tip slide title text -----------------------------------
$slide2 = $p.Slides.Item(2)
$titleShape = $null
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $sh = $slide2.Shapes.Item($i)
    if ($sh.Name -eq "Title 1") { $titleShape = $sh }
}

$tr = $titleShape.TextFrame.TextRange
$fullText = $tr.Text

# The two lines are joined by a manual line break (vertical-tab, chr 11)
# that must stay untouched so the <a:br/> element is preserved.
$breakPos = $fullText.IndexOf([char]11)
$firstLen = $breakPos
$tailStart = $breakPos + 2

# First line: "Tip #1" -> "Tip #6" (keeps the <a:br/> and the second run
# untouched).
$tr.Characters(1, $firstLen).Text = "Tip #6"

# Second line (after the manual line break): replace the remainder of
# the text range with the new tip description.
$tr2 = $titleShape.TextFrame.TextRange
$tailLen = $tr2.Length - ($tailStart - 1)
$tr2.Characters($tailStart, $tailLen).Text = "Avoid using too many arguments"

Write-Output "done"
